$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on column A (Beteckning)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C holds "Förändrad" (Changed) date serials.
# Update every data row (2..lastRow) from 45177 to 45178 (one day later).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
